$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.453991889953613
$ws.Range("B1").Value = 3.219793081283569
$ws.Range("C1").Value = 2.713232517242432
$ws.Range("D1").Value = 2.077782392501831
$ws.Range("E1").Value = 1.251464128494263
